# Update the "Data Table" workbook (Sheet1) with refreshed release-date
# information for the Key Stage 4 / Key Stage 5 destinations rows, and
# move the active selection to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: "Key Stage 4 (KS4) destinations" --------------------------
# Latest period (release date) - corrected release date for 19/20 learners
$ws.Range("C11").Value = "Aug 2020 -  Jul 2021 (19/20 learners) (02/02/23)"
# Next period (release date) - new period added, stored as text so the
# value is not reinterpreted as a date
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "Aug 2021 -  Jul 2022 (20/21 learners) (Oct 23)"

# --- Row 12: "Key Stage 5 (KS5) destinations" --------------------------
$ws.Range("C12").Value = "Aug 2020 -  Jul 2021 (19/20 learners) (02/02/23)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "Aug 2021 -  Jul 2022 (20/21 learners) (Oct 23)"

# --- View state ----------------------------------------------------------
# Reflect the author's final on-screen selection (column D, rows 11-12)
$ws.Activate()
$ws.Range("D11:D12").Select()
